$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

# --- Overview sheet: row 7 is the 826b647f-... file; zh-cn/de-de status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E7").Value = $newStatus
$wsOverview.Range("F7").Value = $newStatus

# --- zh-cn sheet: Status column (C) + Error Detail column (P) for row 7 ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C7").Value = $newStatus
$wsZh.Range("P7").Value = "Handback file name: djsaexry.m4b is different with handoff file name: 826b647f-7625-49ac-bd6a-4a32e34ae767.5eca30f9f16f4457a93b56fe47ccf74bd782c2d6.zh-cn."
# Excel pads character-width input with cell-margin pixels when storing the
# column width; 39.14 is the character width that round-trips to a stored
# width of exactly 40 (matches the target OOXML <col .../> width).
$wsZh.Columns.Item(16).ColumnWidth = 39.14

# --- de-de sheet: Status column (C) + Error Detail column (P) for row 7 ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C7").Value = $newStatus
$wsDe.Range("P7").Value = "Handback file name: djsaexry.m4b is different with handoff file name: 826b647f-7625-49ac-bd6a-4a32e34ae767.5eca30f9f16f4457a93b56fe47ccf74bd782c2d6.de-de."
$wsDe.Columns.Item(16).ColumnWidth = 39.14
